# "Attendance Module is optimized"
#
# Changes applied:
#  1. TestSuite!B4 (Runmode for TC001_VerifyImagerequired): "N" -> "Y"
#  2. TC002_VerifyLogin!B11 / C11: "arjundel" -> "vishaldel"
#     (this removes the only use of "arjundel", so it drops out of the
#     shared-string table and every later string shifts down by one
#     index automatically on save)
#  3. View-state tweaks: per-sheet selection, TestSuite zoom 85 -> 70,
#     and the active sheet/tab moves from TC002_VerifyLogin to
#     TC001_VerifyImagerequired.

$wb = $excel.ActiveWorkbook

$wsTestSuite               = $wb.Worksheets.Item(1)
$wsVerifyAppVersion        = $wb.Worksheets.Item(2)
$wsVerifyLogin             = $wb.Worksheets.Item(3)
$wsVerifyImagerequired     = $wb.Worksheets.Item(4)
$wsVerifyAttendance        = $wb.Worksheets.Item(5)
$wsVerifyDownloadcalls     = $wb.Worksheets.Item(6)
$wsVerifyResourcecentre    = $wb.Worksheets.Item(7)
$wsVerifyVistorlogin       = $wb.Worksheets.Item(8)

# --- data edits -----------------------------------------------------

# Runmode for TC001_VerifyImagerequired flips from N to Y
$wsTestSuite.Range("B4").Value = "Y"

# The login test data row that used to read "arjundel" now reuses
# "vishaldel" (same value already used in B2/B3 of this sheet) --
# "arjundel" is no longer referenced anywhere in the workbook.
$wsVerifyLogin.Range("B11").Value = "vishaldel"
$wsVerifyLogin.Range("C11").Value = "vishaldel"

# --- view-state edits (selection / zoom / active tab) ---------------

# Sheets that are NOT becoming the final active tab: just move their
# saved selection (and, for TestSuite, the zoom level).
[void]$wsTestSuite.Activate()
$excel.ActiveWindow.Zoom = 70
[void]$wsTestSuite.Range("E8").Select()

[void]$wsVerifyLogin.Activate()
[void]$wsVerifyLogin.Range("C11").Select()

[void]$wsVerifyAttendance.Activate()
[void]$wsVerifyAttendance.Range("F13").Select()

# TC001_VerifyImagerequired becomes the active sheet (workbook
# activeTab moves to it, and it picks up tabSelected="1"), with its
# selection moved to H10.
[void]$wsVerifyImagerequired.Activate()
[void]$wsVerifyImagerequired.Range("H10").Select()
